# Update Betfair Back/Lay odds for 2026-02-12 workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (Sion vs FC Basel)
$ws.Range("F8").Value = 2.68
$ws.Range("I8").Value = 2.76
$ws.Range("P8").Value = 2.24

# Row 9 (Thun vs Lausanne)
$ws.Range("F9").Value = 2.08
$ws.Range("G9").Value = 2.22
$ws.Range("H9").Value = 3.85
$ws.Range("I9").Value = 4.4
$ws.Range("J9").Value = 3.35
$ws.Range("P9").Value = 2.2
$ws.Range("Q9").Value = 1.72

# Row 10 (Brentford vs Arsenal)
$ws.Range("H10").Value = 1.75
$ws.Range("I10").Value = 1.76
$ws.Range("AJ10").Value = 140
$ws.Range("AL10").Value = 90
$ws.Range("AM10").Value = 130

# Row 11 (Boyaca Chico vs Jaguares de Cordoba)
$ws.Range("F11").Value = 1.97
$ws.Range("G11").Value = 2.64
$ws.Range("H11").Value = 3.4
$ws.Range("J11").Value = 2.72
$ws.Range("K11").Value = 5.1
$ws.Range("P11").Value = 1.42
$ws.Range("Q11").Value = 2.46

# Row 12 (Athletico-PR vs Santos)
$ws.Range("F12").Value = 2.26
$ws.Range("H12").Value = 3.7
$ws.Range("Q12").Value = 2.22

# Row 13 (Fluminense vs Botafogo FR)
$ws.Range("H13").Value = 3.8

# Row 14 (Corinthians vs Red Bull Bragantino)
$ws.Range("I14").Value = 4.8
$ws.Range("J14").Value = 3.05
$ws.Range("P14").Value = 1.69

# Row 16 (Once Caldas vs Junior FC Barranquilla)
$ws.Range("F16").Value = 2.2
$ws.Range("G16").Value = 2.68
$ws.Range("H16").Value = 3.45
$ws.Range("I16").Value = 4.4
$ws.Range("J16").Value = 2.9
$ws.Range("K16").Value = 3.65
$ws.Range("P16").Value = 1.78
$ws.Range("Q16").Value = 2.08
